$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" column header in H1, copying the header style from G1 ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# --- Updated D/E values for the 100-iteration block (rows 2-11), refit results ---
$ws.Range("D2").Value = 0.5499249211086431
$ws.Range("E2").Value = 0.5499249211086431

$ws.Range("D3").Value = 0.5797913003212973
$ws.Range("E3").Value = 0.5797913003212973

$ws.Range("D4").Value = 0.5838625200862263
$ws.Range("E4").Value = 0.5838625200862263

$ws.Range("D5").Value = 0.5192120555754074
$ws.Range("E5").Value = 0.5192120555754074

$ws.Range("D6").Value = 0.6195488648409917
$ws.Range("E6").Value = 0.6195488648409917

$ws.Range("D7").Value = 0.5281095856855432
$ws.Range("E7").Value = 0.4718904143144568

$ws.Range("D8").Value = 0.5367856328891724
$ws.Range("E8").Value = 0.4632143671108276

$ws.Range("D9").Value = 0.5731840375899194
$ws.Range("E9").Value = 0.4268159624100806

$ws.Range("D10").Value = 0.5330237361742666
$ws.Range("E10").Value = 0.4669762638257334

$ws.Range("D11").Value = 0.5238122463399579
$ws.Range("E11").Value = 0.4761877536600421

# --- New "Label" values (H column): 0 for Control rows, 1 for MDD rows ---
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
